$d = $word.ActiveDocument

# Find the "Kehidupan Sehari-hari" paragraph and the "Akhlak Mulia" paragraph
# that directly follows it (new chapter 29 content), then:
#  1. Apply yellow highlighting to "Kehidupan Sehari-hari" (including the
#     paragraph mark, so the run list + pPr all get <w:highlight w:val="yellow"/>)
#  2. Move the "_GoBack" bookmark from the (now-empty) paragraph after
#     "Akhlak Mulia" to the start of the "Akhlak Mulia" paragraph itself.

$targetIndex = $Null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Kehidupan Sehari-hari")) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq $Null) {
    throw "Could not find 'Kehidupan Sehari-hari' paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)

# 1. Highlight the whole paragraph (runs + paragraph mark) yellow.
$target.Range.Font.HighlightColorIndex = 7

# 2. Move the _GoBack bookmark to the start of the next paragraph
#    ("Akhlak Mulia"), removing it from wherever it currently sits.
$akhlakMulia = $d.Paragraphs.Item($targetIndex + 1)
$bmRange = $d.Range($akhlakMulia.Range.Start, $akhlakMulia.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
